# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    3  = @{ B = 0.3048080303191223; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 2.626907116734944 }
    4  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    6  = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 26.21740644021617;  E = 0.496779210170732; G = 27.15485763464676 }
    7  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 4.429675500412797 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    9  = @{ B = 0.003994804209775715; C = 0.04240448674262143; D = 0.1575252929769615; E = 0.496779210170732; G = 0.7007037941000906 }
    10 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    11 = @{ B = 1.459612070389937;  C = 0.3127903958511391; D = 3.900430680208489;  E = 0.496779210170732; G = 6.169612356620297 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
